# chore: adapt column header formatting to respective input file names
#
# Renames the "_old"/"_new" suffixes used in the header row to the
# concrete format-version identifiers ("_FV2404"/"_FV2410"), turns the
# data range into an Excel Table (ListObject) so the headers/autofilter
# are backed by a proper table definition, and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the data (header row + all data rows/cols).
$used = $ws.UsedRange
$lastRow = $used.Rows.Count
$lastCol = $used.Columns.Count

# --- 1. Rename header cells: "<Name>_old" -> "<Name>_FV2404",
#        "<Name>_new" -> "<Name>_FV2410" -------------------------------
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2404"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2410"
        }
    }
}

# --- 2. Turn the header+data range into an Excel Table -----------------
$rng = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row ------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
